$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "Mon&Tues" remaining-time column (D3:D17) with updated values
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1.5
$ws.Range("D11").Value = 0.5
$ws.Range("D12").Value = 0
$ws.Range("D13").Value = 0.5
$ws.Range("D14").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("D17").Value = 0

# Update the active selection/view to match author's saved state
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("E18").Select()
